# Remove the blank leading row above the header on the
# "propublica_odd_comparison" sheet, shifting all data up by one row
# (B3:AB33 -> B2:AB32).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("propublica_odd_comparison")
$ws.Activate()

$ws.Rows.Item(1).Delete()

$ws.Rows.Item(2).Select()
